$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data: turn row 1 into a header row, shift the existing
#     5 data rows down to rows 2-6, and add a new "Formatted Loyalty" column ---

$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "Loyalty Number"
$ws.Range("B1").Value = "Total Transactions"
$ws.Range("C1").Value = "Date of Issue"
$ws.Range("D1").Value = "Formatted Loyalty"

# New computed column: loyalty number + formatted issue date
$ws.Range("D2:D6").Formula = "=A2&`" `"&TEXT(C2, `"MM/DD/yyyy`")"

# Column D (and the C1 header) use the same date-style look (numFmt 14,
# no wrap) that the rest of column C already had.
$ws.Range("D2:D6").NumberFormat = "m/d/yyyy"
$ws.Range("C1").NumberFormat = "m/d/yyyy"

# The first data row's Date of Issue cell (C2) gets its own distinct
# number format.
$ws.Range("C2").NumberFormat = "m/d/yyyy;@"

# Cursor / print setup bookkeeping to match the resaved workbook.
$ws.Range("E8").Select()
$ws.PageSetup.Orientation = 1

# --- Pie chart: retarget it from "Loyalty Number" to "Store" ---

$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Top 5 Total Transactions by Store"

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Sheet1'!B1,'Sheet1'!`$D`$2:`$D`$6,'Sheet1'!`$B`$2:`$B`$6,1)"
